# Updates cryptocurrency price/volume data in the worksheet to match the
# latest scrape results (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.380.22'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').Value = '3.601.00'
$ws.Range('E3').Value = '  -0.84%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.79'
$ws.Range('E5').Value = '  -2.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '189.03'
$ws.Range('E6').Value = '  -3.04%  '
$ws.Range('B7').Value = 'LidoStakedEther'
$ws.Range('C7').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D7').Value = '3.596.15'
$ws.Range('E7').Value = '  -0.78%  '
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.630'
$ws.Range('E8').Value = '  -2.42%  '
$ws.Range('E9').Value = '  +0.08%  '
$ws.Range('E10').Value = '  +4.04%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.661'
$ws.Range('E11').Value = '  -1.36%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '55.90'
$ws.Range('E12').Value = '  -4.63%  '
$ws.Range('E13').Value = '  +7.28%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.67'
$ws.Range('E14').Value = '  -2.68%  '
$ws.Range('D15').Value = '4.179.64'
$ws.Range('E15').Value = '  -0.76%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '19.79'
$ws.Range('E16').Value = '  -0.06%  '
$ws.Range('D17').Value = '3.602.86'
$ws.Range('E17').Value = '  -0.80%  '
$ws.Range('D18').Value = '70.316.16'
$ws.Range('E18').Value = '  -0.12%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.63'
$ws.Range('E19').Value = '  -0.87%  '
$ws.Range('E20').Value = '  +0.05%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.04'
$ws.Range('E21').Value = '  -2.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '489.11'
$ws.Range('E22').Value = '  +0.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '19.46'
$ws.Range('E23').Value = '  +1.30%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.85'
$ws.Range('E24').Value = '  -9.18%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '96.52'
$ws.Range('E25').Value = '  +5.75%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.35'
$ws.Range('E26').Value = '  -2.56%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.97'
$ws.Range('E27').Value = '  -6.20%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.00'
$ws.Range('E28').Value = '  -4.70%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.36'
$ws.Range('E29').Value = '  -3.22%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.21'
$ws.Range('E30').Value = '  -2.25%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.61'
$ws.Range('E31').Value = '  -4.19%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.21'
$ws.Range('E32').Value = '  -0.62%  '
$ws.Range('E33').Value = '  -3.24%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '65.74'
$ws.Range('E34').Value = '  -0.33%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '571.87'
$ws.Range('E35').Value = '  -8.64%  '
$ws.Range('E36').Value = '  -4.30%  '
$ws.Range('D37').Value = '0.0₃0809'
$ws.Range('E37').Value = '  -1.91%  '
$ws.Range('E38').Value = '  +0.08%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.396'
$ws.Range('E39').Value = '  -4.60%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.31'
$ws.Range('E40').Value = '  +14.64%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.97'
$ws.Range('E41').Value = '  +4.83%  '
$ws.Range('E42').Value = '  -2.94%  '
$ws.Range('E43').Value = '  -6.25%  '
$ws.Range('B44').Value = 'ApeXProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.55'
$ws.Range('E44').Value = '  +6.92%  '
$ws.Range('B45').Value = 'ThetaToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.02'
$ws.Range('E45').Value = '  -4.30%  '
$ws.Range('D46').Value = '3.213.46'
$ws.Range('E46').Value = '  -2.38%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0444'
$ws.Range('E47').Value = '  -2.06%  '
$ws.Range('E48').Value = '  +4.77%  '
$ws.Range('E49').Value = '  -0.53%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.00'
$ws.Range('E50').Value = '  -0.04%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.22'
$ws.Range('E51').Value = '  -4.04%  '
